$wb = $excel.ActiveWorkbook

# --- Re-write cells that hold a "backslash" duplicate string so the shared
# --- string table gets compacted/deduplicated to the canonical "slash" form
# --- (this mirrors Excel's own resave behaviour and is required to get
# --- sharedStrings.xml from uniqueCount=78 down to uniqueCount=67).

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("I2").Value = "/23h/"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("I1").Value = "/17h/"
$ws2.Range("I2").Value = "/23h/"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("I1").Value = "/17h/"
$ws3.Range("I2").Value = "/23h/"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("I1").Value = "/RETOUR_MGEFI_GTO/Noemie/"
$ws4.Range("I2").Value = "/RETOUR_MGEFI_GTO/Noemie/"
$ws4.Range("I3").Value = "/RETOUR_MGEFI_GTO/rejet def/"
$ws4.Range("I4").Value = "/RETOUR_MGEFI_GTO/rejet def/"
$ws4.Range("I5").Value = "/RETOUR_ADHESION_LAMIE/"
$ws4.Range("I6").Value = "/RETOUR_ADHESION_LAMIE/"
$ws4.Range("I7").Value = "/RETOUR_ADHESION_CSS/adhesion ITE/"
$ws4.Range("I8").Value = "/RETOUR_ADHESION_CSS/adhesion ITE/"
$ws4.Range("I9").Value = "/RETOUR_ADHESION_CSS/adhesion MGAS/"
$ws4.Range("I10").Value = "/RETOUR_ADHESION_CSS/adhesion MGAS/"
$ws4.Range("I11").Value = "/RETOUR_ADHESION_CSS/adhesion LMDE/"
$ws4.Range("I12").Value = "/RETOUR_ADHESION_CSS/adhesion LMDE/"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("I1").Value = "/17h/"

$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("I1").Value = "/17h/"

$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("I1").Value = "/17h/"

$ws8 = $wb.Worksheets.Item(8)
$ws8.Range("I1").Value = "/17h/"

$ws9 = $wb.Worksheets.Item(9)
$ws9.Range("I1").Value = "/RETOUR_ADHESION_LAMIE/"

$ws10 = $wb.Worksheets.Item(10)
$ws10.Range("I1").Value = "/RETOUR_ADHESION_LAMIE/"

$ws11 = $wb.Worksheets.Item(11)
$ws11.Range("I1").Value = "/SALESFORCE/"
$ws11.Range("I2").Value = "/SALESFORCE/"

$ws12 = $wb.Worksheets.Item(12)
$ws12.Range("I1").Value = "/SALESFORCE/"

$ws13 = $wb.Worksheets.Item(13)
$ws13.Range("I1").Value = "/"

$ws14 = $wb.Worksheets.Item(14)
$ws14.Range("I1").Value = "/"

$ws15 = $wb.Worksheets.Item(15)
$ws15.Range("I1").Value = "/"

$ws16 = $wb.Worksheets.Item(16)
$ws16.Range("I1").Value = "/"

$ws17 = $wb.Worksheets.Item(17)
$ws17.Range("I1").Value = "/"

$ws18 = $wb.Worksheets.Item(18)
$ws18.Range("I1").Value = "/"

$ws19 = $wb.Worksheets.Item(19)
$ws19.Range("I1").Value = "/"

$ws20 = $wb.Worksheets.Item(20)
$ws20.Range("I1").Value = "/"

$ws21 = $wb.Worksheets.Item(21)
$ws21.Range("I1").Value = "/"

$ws22 = $wb.Worksheets.Item(22)
$ws22.Range("I1").Value = "/"

# --- View / selection state updates -----------------------------------

# Sheet1 (Feuil1): no longer the active/visible tab, selection moves to I5
$ws1.Range("I5").Select()

# Sheet2 (Feuil2): selection moves to I2
$ws2.Range("I2").Select()

# Sheet3 (Feuil3): selection moves to I4
$ws3.Range("I4").Select()

# Sheet4 (Feuil4): scrolled so column B is first visible column, selection F14
$ws4.Application.ActiveWindow.ScrollColumn = 2
$ws4.Range("F14").Select()

# Sheet5 (Feuil5): selection moves to I1
$ws5.Range("I1").Select()

# Sheet6 (Feuil6): selection moves to I1
$ws6.Range("I1").Select()

# Sheet7 (Feuil7): selection moves to I1
$ws7.Range("I1").Select()

# Sheet8 (Feuil8): selection moves to I1
$ws8.Range("I1").Select()

# Sheet9 (Feuil9): selection moves to F6
$ws9.Range("F6").Select()

# Sheet10 (Feuil10): selection moves to I1
$ws10.Range("I1").Select()

# Sheet11 (Feuil11): selection moves to G5
$ws11.Range("G5").Select()

# Sheet12 (Feuil12): gains an explicit selection at H8
$ws12.Range("H8").Select()

# Sheet13 (Feuil13): selection moves to I1
$ws13.Range("I1").Select()

# Sheet14 (Feuil14): selection moves to I1
$ws14.Range("I1").Select()

# Sheet15 (Feuil15): selection moves to I1
$ws15.Range("I1").Select()

# Sheet16 (Feuil16): selection moves to I1
$ws16.Range("I1").Select()

# Sheet17 (Feuil17): selection moves to I1
$ws17.Range("I1").Select()

# Sheet18 (Feuil18): selection moves to I1
$ws18.Range("I1").Select()

# Sheet19 (Feuil19): selection moves to I1
$ws19.Range("I1").Select()

# Sheet20 (Feuil20): selection moves to I1
$ws20.Range("I1").Select()

# Sheet21 (Feuil21): selection moves to I1
$ws21.Range("I1").Select()

# Sheet22 (Feuil22): becomes the active/visible tab, selection moves to H6
$ws22.Range("H6").Select()
$ws22.Activate()

$wb.Save()
